# Applies updated cryptocurrency price/volume/date/hour data to Sheet1
# (rows 2-51, columns D:Price, E:Volume(1h), F:Data, G:Hora)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row number, then a hashtable of column letter -> new text value
$updates = @(
    @{ Row=2; Cols=@{ D='326.73'; E='2.41%'; F='3-2-2023'; G='1' } },
    @{ Row=3; Cols=@{ D='39.59'; E='0.12%'; F='3-2-2023'; G='1' } },
    @{ Row=4; Cols=@{ D='5.879'; E='14.35%'; F='3-2-2023'; G='1' } },
    @{ Row=5; Cols=@{ D='0.08026'; E='-1.89%'; F='3-2-2023'; G='1' } },
    @{ Row=6; Cols=@{ D='4.582'; E='2.44%'; F='3-2-2023'; G='1' } },
    @{ Row=7; Cols=@{ D='8.696'; E='1.42%'; F='3-2-2023'; G='1' } },
    @{ Row=8; Cols=@{ D='1.912'; E='-0.30%'; F='3-2-2023'; G='1' } },
    @{ Row=9; Cols=@{ D='2.940'; E='-0.17%'; F='3-2-2023'; G='1' } },
    @{ Row=10; Cols=@{ D='0.9331'; E='-1.20%'; F='3-2-2023'; G='1' } },
    @{ Row=11; Cols=@{ D='0.1240'; E='-3.64%'; F='3-2-2023'; G='1' } },
    @{ Row=12; Cols=@{ D='0.1965'; E='0.37%'; F='3-2-2023'; G='1' } },
    @{ Row=13; Cols=@{ D='8.759'; E='30.47%'; F='3-2-2023'; G='1' } },
    @{ Row=14; Cols=@{ D='0.09136'; E='-0.24%'; F='3-2-2023'; G='1' } },
    @{ Row=15; Cols=@{ D='0.03521'; E='2.17%'; F='3-2-2023'; G='1' } },
    @{ Row=16; Cols=@{ D='0.09628'; E='1.12%'; F='3-2-2023'; G='1' } },
    @{ Row=17; Cols=@{ D='0.001301'; E='-7.57%'; F='3-2-2023'; G='1' } },
    @{ Row=18; Cols=@{ D='0.006148'; E='4.78%'; F='3-2-2023'; G='1' } },
    @{ Row=19; Cols=@{ D='3.338'; E='-0.79%'; F='3-2-2023'; G='1' } },
    @{ Row=20; Cols=@{ D='0.3537'; E='-0.05%'; F='3-2-2023'; G='1' } },
    @{ Row=21; Cols=@{ D='0.1431'; E='8.41%'; F='3-2-2023'; G='1' } },
    @{ Row=22; Cols=@{ D='0.2415'; E='5.50%'; F='3-2-2023'; G='1' } },
    @{ Row=23; Cols=@{ D='0.04433'; E='1.13%'; F='3-2-2023'; G='1' } },
    @{ Row=24; Cols=@{ E='3.34%'; F='3-2-2023'; G='1' } },
    @{ Row=25; Cols=@{ D='0.004348'; E='-0.60%'; F='3-2-2023'; G='1' } },
    @{ Row=26; Cols=@{ D='0.0001142'; E='0.71%'; F='3-2-2023'; G='1' } },
    @{ Row=27; Cols=@{ E='1.08%'; F='3-2-2023'; G='1' } },
    @{ Row=28; Cols=@{ F='3-2-2023'; G='1' } },
    @{ Row=29; Cols=@{ F='3-2-2023'; G='1' } },
    @{ Row=30; Cols=@{ F='3-2-2023'; G='1' } },
    @{ Row=31; Cols=@{ F='3-2-2023'; G='1' } },
    @{ Row=32; Cols=@{ F='3-2-2023'; G='1' } },
    @{ Row=33; Cols=@{ F='3-2-2023'; G='1' } },
    @{ Row=34; Cols=@{ F='3-2-2023'; G='1' } },
    @{ Row=35; Cols=@{ F='3-2-2023'; G='1' } },
    @{ Row=36; Cols=@{ F='3-2-2023'; G='1' } },
    @{ Row=37; Cols=@{ F='3-2-2023'; G='1' } },
    @{ Row=38; Cols=@{ F='3-2-2023'; G='1' } },
    @{ Row=39; Cols=@{ D='0.02422'; E='-0.92%'; F='3-2-2023'; G='1' } },
    @{ Row=40; Cols=@{ D='0.05209'; E='0.10%'; F='3-2-2023'; G='1' } },
    @{ Row=41; Cols=@{ D='0.007443'; E='-3.31%'; F='3-2-2023'; G='1' } },
    @{ Row=42; Cols=@{ D='0.1407'; E='-2.14%'; F='3-2-2023'; G='1' } },
    @{ Row=43; Cols=@{ D='0.008699'; E='0.80%'; F='3-2-2023'; G='1' } },
    @{ Row=44; Cols=@{ D='0.002123'; E='4.21%'; F='3-2-2023'; G='1' } },
    @{ Row=45; Cols=@{ D='0.009638'; E='4.38%'; F='3-2-2023'; G='1' } },
    @{ Row=46; Cols=@{ D='0.00006691'; E='3.77%'; F='3-2-2023'; G='1' } },
    @{ Row=47; Cols=@{ D='0.00000000751'; E='0.69%'; F='3-2-2023'; G='1' } },
    @{ Row=48; Cols=@{ D='0.003003'; E='5.62%'; F='3-2-2023'; G='1' } },
    @{ Row=49; Cols=@{ D='0.001423'; E='-42.30%'; F='3-2-2023'; G='1' } },
    @{ Row=50; Cols=@{ D='0.00002103'; E='0.69%'; F='3-2-2023'; G='1' } },
    @{ Row=51; Cols=@{ D='0.0002003'; E='0.69%'; F='3-2-2023'; G='1' } }
)

foreach ($update in $updates) {
    $r = $update.Row
    foreach ($col in $update.Cols.Keys) {
        $cell = $ws.Range("$col$r")
        # Force text storage so values like "326.73" or "3-2-2023" are not
        # reinterpreted by Excel as a number/date, matching the inlineStr cells
        $cell.NumberFormat = "@"
        $cell.Value = $update.Cols[$col]
        # Reset style back to Normal so no stray number-format style lingers
        $cell.Style = "Normal"
    }
}
